# Update "想去人数" (interested-attendee count) figures to match the
# latest scrape results reflected in the regenerated gh-pages output.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 292
$wsExhibit.Range("F4").Value = 1120
$wsExhibit.Range("F5").Value = 583

# Sheet "全部类型" (all types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 292
$wsAll.Range("F4").Value = 1120
$wsAll.Range("F6").Value = 583
